# test(web)/qa/salesforce/residencial: agregar planes de residencial 2p BRM
#
# Fills in the "Plans" sheet with additional Residencial "2p" (Sin_TotalPlay_TV)
# plan rows, switching plan #1/#2's service type to Sin_TotalPlay_TV and adding
# three new plan rows (#3, #4, #5) with Residencial / Sin_TotalPlay_TV plans at
# 200, 500 and 1000 Megas respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plans")

# Plan #1 (row 5): switch service type to Sin_TotalPlay_TV and its Megas to 50
$ws.Range("C5").Value = "Sin_TotalPlay_TV"
$ws.Range("D5").Value = 50

# Plan #2 (row 6): switch service type to Sin_TotalPlay_TV (Megas unchanged, 100)
$ws.Range("C6").Value = "Sin_TotalPlay_TV"

# Plan #3 (row 7): new Residencial / Sin_TotalPlay_TV plan at 200 Megas
$ws.Range("B7").Value = "Residencial"
$ws.Range("C7").Value = "Sin_TotalPlay_TV"
$ws.Range("D7").Value = 200

# Plan #4 (row 8): new Residencial / Sin_TotalPlay_TV plan at 500 Megas
$ws.Range("B8").Value = "Residencial"
$ws.Range("C8").Value = "Sin_TotalPlay_TV"
$ws.Range("D8").Value = 500

# Plan #5 (row 9): new Residencial / Sin_TotalPlay_TV plan at 1000 Megas
$ws.Range("B9").Value = "Residencial"
$ws.Range("C9").Value = "Sin_TotalPlay_TV"
$ws.Range("D9").Value = 1000

# Leave the cursor where the author last left it while editing
$ws.Activate()
$ws.Range("E11").Select() | Out-Null
